# Update message design: add "Matching Game" sheet documenting CMD_GAME_MATCHING,
# move the active/selected tab from "Login,Register" to "Chat", and adjust the
# remembered selections on a couple of sheets.

$wb = $excel.ActiveWorkbook

$wsLogin  = $wb.Worksheets.Item("Login,Register")
$wsChat   = $wb.Worksheets.Item("Chat")

# --- New worksheet: "Matching Game" -----------------------------------------
# Added after the last existing sheet ("Chat") so it lands at the end of the
# tab strip.
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Matching Game"

# Row 1 - header row (reuse formatting from the "Login,Register" sheet header)
$wsLogin.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Command ID"

$wsLogin.Range("C1:E1").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Field"
$ws.Range("C1").Value = "ValueType"
$ws.Range("D1").Value = "Description"

$wsLogin.Range("K1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)

# Row 2 - CMD_GAME_MATCHING / code / int(0,1) / 0:fail,1:success
$wsLogin.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "CMD_GAME_MATCHING"

$wsLogin.Range("C2:D2").Copy()
$ws.Range("B2:C2").PasteSpecial(-4122)
$ws.Range("B2").Value = "code"
$ws.Range("C2").Value = "int(0,1)"

$wsLogin.Range("C2").Copy()
$ws.Range("D2:G2").PasteSpecial(-4122)
$ws.Range("D2").Value = "0:fail,1:success"

# Row 3 - room_id / int / the room id user joined
$wsLogin.Range("C2").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)
$ws.Range("B3").Value = "room_id"
$ws.Range("C3").Value = "int"
$ws.Range("D3").Value = "the room id user joined"

# Rows 4-5 - blank, formatted body rows
$wsLogin.Range("C2").Copy()
$ws.Range("A4:G5").PasteSpecial(-4122)

# Row 6 - trailing note, unstyled
$ws.Range("A6").Value = "Sau khi nhan message nay, player se nhan duoc message CMD_FRIEND_INFO"

$ws.Range("A12").Select()

# --- Selection bookkeeping to match the new active tab ----------------------
# "Login,Register" loses its remembered selection / tab flag, "Chat" gains it.
$wsLogin.Range("B1:H5").Select()

$wsChat.Activate()
$wsChat.Range("F11").Select()

$wb.Application.ScreenUpdating = $true
